$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 40

# Date/Week columns get Excel "smart" type inference (date, number) on plain
# Value assignment; force them to Text first, write the value, then drop the
# number-format override so the cell ends up as plain, unstyled text.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2024-01-09"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "22:54:46"

$ws.Cells.Item($row, 3).Value = "Tuesday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "01"
$ws.Cells.Item($row, 4).ClearFormats()

$ws.Cells.Item($row, 5).Value = 139477
$ws.Cells.Item($row, 6).Value = 142691
$ws.Cells.Item($row, 7).Value = 172105
$ws.Cells.Item($row, 8).Value = 147596
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 118510
$ws.Cells.Item($row, 11).Value = 224742
$ws.Cells.Item($row, 12).Value = 250823
$ws.Cells.Item($row, 13).Value = 185112
$ws.Cells.Item($row, 14).Value = 110455
$ws.Cells.Item($row, 15).Value = 40738
$ws.Cells.Item($row, 16).Value = 30875
$ws.Cells.Item($row, 17).Value = 72616
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42253
$ws.Cells.Item($row, 20).Value = -1
